$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column D (Price) values are free-form text (e.g. "58.336.32", "0.999")
# that Excel would otherwise auto-coerce into numbers and mangle (trailing
# zeros dropped, floating point noise). Force text format first so the
# literal string is preserved exactly, matching the inlineStr cells in the file.
$dCells = 'D2', 'D3', 'D5', 'D6', 'D8', 'D9', 'D12', 'D13', 'D14', 'D15', 'D17', 'D18', 'D20', 'D21', 'D22', 'D23', 'D24', 'D27', 'D29', 'D31', 'D34', 'D35', 'D37', 'D38', 'D39', 'D40', 'D41', 'D42', 'D43', 'D44', 'D45', 'D46', 'D47', 'D48', 'D49', 'D50', 'D51'
foreach ($addr in $dCells) {
    $ws.Range($addr).NumberFormat = "@"
}

$ws.Range('D2').Value = '58.071.82'
$ws.Range('E2').Value = '  -2.02%  '
$ws.Range('D3').Value = '2.473.79'
$ws.Range('E3').Value = '  -2.14%  '
$ws.Range('E4').Value = '  +0.21%  '
$ws.Range('D5').Value = '522.85'
$ws.Range('E5').Value = '  -3.13%  '
$ws.Range('D6').Value = '132.68'
$ws.Range('E6').Value = '  -4.01%  '
$ws.Range('E7').Value = '  +0.12%  '
$ws.Range('D8').Value = '0.559'
$ws.Range('E8').Value = '  -1.51%  '
$ws.Range('D9').Value = '0.0997'
$ws.Range('E9').Value = '  -2.02%  '
$ws.Range('E11').Value = '  +0.75%  '
$ws.Range('D12').Value = '0.342'
$ws.Range('E12').Value = '  -2.19%  '
$ws.Range('D13').Value = '2.920.22'
$ws.Range('E13').Value = '  -1.87%  '
$ws.Range('D14').Value = '58.176.71'
$ws.Range('E14').Value = '  -1.75%  '
$ws.Range('D15').Value = '22.30'
$ws.Range('E15').Value = '  -4.04%  '
$ws.Range('E16').Value = '  -2.30%  '
$ws.Range('D17').Value = '2.481.40'
$ws.Range('E17').Value = '  -2.04%  '
$ws.Range('D18').Value = '10.90'
$ws.Range('E18').Value = '  -2.18%  '
$ws.Range('E19').Value = '  -2.66%  '
$ws.Range('D20').Value = '321.13'
$ws.Range('E20').Value = '  -1.64%  '
$ws.Range('D21').Value = '0.999'
$ws.Range('D22').Value = '5.79'
$ws.Range('E22').Value = '  -3.02%  '
$ws.Range('D23').Value = '64.30'
$ws.Range('E23').Value = '  -1.81%  '
$ws.Range('D24').Value = '0.411'
$ws.Range('E24').Value = '  -3.46%  '
$ws.Range('E25').Value = '  -0.40%  '
$ws.Range('E26').Value = '  -3.42%  '
$ws.Range('D27').Value = '7.46'
$ws.Range('E27').Value = '  -2.78%  '
$ws.Range('E28').Value = '  -3.71%  '
$ws.Range('D29').Value = '6.39'
$ws.Range('E29').Value = '  -5.28%  '
$ws.Range('E30').Value = '  -4.54%  '
$ws.Range('D31').Value = '166.33'
$ws.Range('E31').Value = '  -0.77%  '
$ws.Range('E32').Value = '  -4.74%  '
$ws.Range('E33').Value = '  -0.01%  '
$ws.Range('D34').Value = '1.00'
$ws.Range('E34').Value = '  +0.31%  '
$ws.Range('D35').Value = '18.21'
$ws.Range('E35').Value = '  -1.63%  '
$ws.Range('E36').Value = '  -9.03%  '
$ws.Range('D37').Value = '3.97'
$ws.Range('E37').Value = '  -4.19%  '
$ws.Range('D38').Value = '1.49'
$ws.Range('E38').Value = '  -4.63%  '
$ws.Range('D39').Value = '0.796'
$ws.Range('E39').Value = '  -3.33%  '
$ws.Range('D40').Value = '3.51'
$ws.Range('E40').Value = '  -3.96%  '
$ws.Range('D41').Value = '276.41'
$ws.Range('E41').Value = '  -3.81%  '
$ws.Range('D42').Value = '4.98'
$ws.Range('E42').Value = '  -4.93%  '
$ws.Range('D43').Value = '0.592'
$ws.Range('E43').Value = '  -2.91%  '
$ws.Range('D44').Value = '127.22'
$ws.Range('E44').Value = '  -3.58%  '
$ws.Range('D45').Value = '0.0911'
$ws.Range('E45').Value = '  -2.48%  '
$ws.Range('D46').Value = '0.0494'
$ws.Range('E46').Value = '  -3.50%  '
$ws.Range('D47').Value = '0.0216'
$ws.Range('E47').Value = '  -2.84%  '
$ws.Range('D48').Value = '17.13'
$ws.Range('E48').Value = '  -1.95%  '
$ws.Range('D49').Value = '1.739.76'
$ws.Range('E49').Value = '  -1.44%  '
$ws.Range('D50').Value = '0.974'
$ws.Range('E50').Value = '  -1.50%  '
$ws.Range('D51').Value = '4.64'
$ws.Range('E51').Value = '  -2.58%  '
